# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.261.08'
$ws.Range('E2').Value = '  +0.61%  '
$ws.Range('D3').Value = '2.483.62'
$ws.Range('E3').Value = '  +1.15%  '
$ws.Range('E4').Value = '  +0.01%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '585.86'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +1.23%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '172.12'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +3.92%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +0.82%  '
$ws.Range('D9').Value = '2.483.42'
$ws.Range('E9').Value = '  +1.22%  '
$ws.Range('E10').Value = '  +3.47%  '
$ws.Range('E11').Value = '  +1.32%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '4.95'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +1.52%  '
$ws.Range('E13').Value = '  +0.59%  '
$ws.Range('E14').Value = '  +1.26%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '25.58'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +1.17%  '
$ws.Range('D16').Value = '66.942.35'
$ws.Range('E16').Value = '  +0.18%  '
$ws.Range('E17').Value = '  +1.67%  '
$ws.Range('D18').Value = '2.482.21'
$ws.Range('E18').Value = '  +1.11%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '7.77'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +0.68%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '11.04'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -2.42%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '352.34'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -0.70%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.00'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -0.40%  '
$ws.Range('E23').Value = '  +0.06%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '69.02'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -0.52%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '4.25'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +0.92%  '
$ws.Range('E26').Value = '  +3.12%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '9.30'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +4.56%  '
$ws.Range('D28').Value = '2.574.80'
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('E29').Value = '  +0.37%  '
$ws.Range('D30').Value = '0.0₃0915'
$ws.Range('E30').Value = '  +1.92%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '511.39'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +1.25%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '7.74'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -0.57%  '
$ws.Range('E33').Value = '  +3.11%  '
$ws.Range('E34').Value = '  +0.13%  '
$ws.Range('E35').Value = '  +0.02%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '162.13'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +1.95%  '
$ws.Range('E37').Value = '  +2.30%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '18.70'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +0.78%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '18.19'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -1.34%  '
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('E42').Value = '  +2.37%  '
$ws.Range('E43').Value = '  +1.30%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '4.85'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +1.88%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '2.39'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +3.27%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '143.98'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +1.95%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.0₆0263'
$ws.Range('E47').Value = '  +4.60%  '
$ws.Range('B48').Value = 'Filecoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '3.50'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +0.79%  '
$ws.Range('E49').Value = '  +0.66%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.0736'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +0.54%  '
$ws.Range('E51').Value = '  +0.04%  '
